$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.902.29"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").Value = "3.408.93"
$ws.Range("E3").Value = "  -1.75%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "575.26"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "147.83"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +1.09%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.94"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.80%  "

$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("E11").Value = "  +2.81%  "

$ws.Range("D12").Value = "3.995.90"
$ws.Range("E12").Value = "  -1.61%  "

$ws.Range("E13").Value = "  -0.02%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "28.39"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.54%  "

$ws.Range("D15").Value = "3.406.22"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "61.936.73"
$ws.Range("E17").Value = "  -1.57%  "

$ws.Range("E18").Value = "  +1.17%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.54"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "8.94"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.12%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "380.07"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.566"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "3.570.06"
$ws.Range("E25").Value = "  -1.06%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0000111"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.00%  "

$ws.Range("E27").Value = "  +0.28%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.61"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.90"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.01%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.12"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  -0.01%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.33"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "23.01"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.85%  "

$ws.Range("E35").Value = "  +3.81%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.61"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("E37").Value = "  -2.17%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "169.53"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("E39").Value = "  -4.97%  "

$ws.Range("D40").Value = "3.445.40"
$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("E41").Value = "  +3.49%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.781"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "42.40"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.36"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.13%  "

$ws.Range("E45").Value = "  -2.49%  "

$ws.Range("E46").Value = "  -3.04%  "

$ws.Range("D47").Value = "2.541.46"
$ws.Range("E47").Value = "  -2.77%  "

$ws.Range("E48").Value = "  +2.53%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "22.69"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.02%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.19"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.88%  "

$ws.Range("E51").Value = "  +0.20%  "
